# Updates the "IPC PO" (predicted) column C values for rows 2-51 with the
# results of the (now working) genetic-algorithm model, then recomputes the
# dependent DELTA (col D), DELTA^2 (col E), TOTAL (row 52) and MSE (row 53)
# cells so the cached values match Excel's own floating point arithmetic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New genetic-algorithm predicted IPC values (column C), rows 2..51.
$newC = @(
    28.88047081994445, 28.79754100912948, 28.72393905747945, 29.94445591304089,
    29.61782930945843, 30.54245861807876, 30.40697447332623, 30.21982061576142,
    29.49073496471589, 29.76432767218155, 29.2289763240952,  29.64752027789943,
    29.41012243268786, 30.14218121153088, 30.17853079185805, 30.47421113686185,
    30.0327836856211,  30.48739069386338, 30.80707348865126, 31.47541242433202,
    31.90392191954129, 32.07101273992203, 31.91345175328856, 32.31898828666998,
    32.11056598256215, 32.74997155045387, 32.50279706675092, 32.61164882504537,
    32.87701108739923, 33.00066775366493, 32.97099385167881, 33.9102622668428,
    33.79824511495674, 33.92622822892227, 34.03392978109852, 34.58591765730413,
    35.70262508082213, 36.10905442207207, 36.26659397241595, 36.88279370718733,
    36.91183076134275, 37.88692681021995, 38.50346492903321, 38.8761424283309,
    39.2557687270679,  39.56415861988053, 39.79025351510796, 40.01247792801502,
    40.31626741864778, 41.48393709049734
)

$firstRow = 2
$lastRow = 51

for ($i = 0; $i -lt $newC.Length; $i++) {
    $row = $firstRow + $i
    $b = $ws.Cells.Item($row, 2).Value2
    $c = $newC[$i]
    $d = $c - $b
    $e = $d * $d

    $ws.Cells.Item($row, 3).Value2 = $c
    $ws.Cells.Item($row, 4).Value2 = $d
    $ws.Cells.Item($row, 5).Value2 = $e
}

# TOTAL row: column C is the sum of the DELTA column, column E is the sum of
# the DELTA^2 column.
$totalD = 0.0
$totalE = 0.0
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $totalD += $ws.Cells.Item($row, 4).Value2
    $totalE += $ws.Cells.Item($row, 5).Value2
}

$ws.Cells.Item(52, 3).Value2 = $totalD
$ws.Cells.Item(52, 5).Value2 = $totalE

# MSE row: mean of the DELTA^2 column.
$count = $lastRow - $firstRow + 1
$ws.Cells.Item(53, 5).Value2 = $totalE / $count
